# Fruta / hortaliza, semanal
# Insert a new weekly record at row 25, pushing the existing rows
# (old 25..103) down to (26..104), and populate the new row with the
# latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25 - this shifts rows 25..103 down to 26..104
# and inherits formatting (e.g. the date style on column D) from the
# surrounding rows automatically.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new weekly record.
$ws.Range("A25").Value = 9
$ws.Range("B25").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C25").Value = "Metropolitana"
$ws.Range("D25").Value = 44600
$ws.Range("E25").Value = 13
$ws.Range("F25").Value = 100112022
$ws.Range("G25").Value = "Arveja Verde"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 28
$ws.Range("K25").Value = 30000
$ws.Range("L25").Value = 30000
$ws.Range("M25").Value = 30000
$ws.Range("N25").Value = "`$/saco 25 kilos"
$ws.Range("O25").Value = "Carahue"
$ws.Range("P25").Value = 1200
$ws.Range("Q25").Value = 25
$ws.Range("R25").Value = "Hortaliza"
